$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 7991
$ws.Range("E2").Value = 523
$ws.Range("F2").Value = 523
$ws.Range("G2").Value = 507
$ws.Range("H2").Value = 366
$ws.Range("I2").Value = 329
$ws.Range("J2").Value = 37
$ws.Range("K2").Value = 64030
$ws.Range("L2").Value = 56794
$ws.Range("M2").Value = 7236
$ws.Range("N2").Value = 6772
$ws.Range("O2").Value = 463
$ws.Range("P2").Value = 2877
$ws.Range("Q2").Value = -2187
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 312
$ws.Range("T2").Value = 22
$ws.Range("V2").Value = 34656
$ws.Range("W2").Value = 6.55
$ws.Range("X2").Value = 4.58
$ws.Range("Y2").Value = 4.93
$ws.Range("Z2").Value = 0.58
$ws.Range("AA2").Value = 784.94
$ws.Range("AB2").Value = 152.16
$ws.Range("AC2").Value = 571
$ws.Range("AD2").Value = 11.27
$ws.Range("AE2").Value = 11827
$ws.Range("AF2").Value = 0.54
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 4.66
$ws.Range("AI2").Value = 44.65
$ws.Range("AJ2").Value = 57545890
$ws.Range("U2").ClearContents()

# --- Row 3 ---
$ws.Range("D3").Value = 7947
$ws.Range("E3").Value = 670
$ws.Range("F3").Value = 670
$ws.Range("G3").Value = 639
$ws.Range("H3").Value = 510
$ws.Range("I3").Value = 461
$ws.Range("J3").Value = 49
$ws.Range("K3").Value = 70439
$ws.Range("L3").Value = 62896
$ws.Range("M3").Value = 7543
$ws.Range("N3").Value = 7076
$ws.Range("O3").Value = 467
$ws.Range("P3").Value = 2877
$ws.Range("Q3").Value = -2879
$ws.Range("R3").Value = 127
$ws.Range("S3").Value = 3592
$ws.Range("T3").Value = 25
$ws.Range("V3").Value = 36675
$ws.Range("W3").Value = 8.43
$ws.Range("X3").Value = 6.42
$ws.Range("Y3").Value = 6.66
$ws.Range("Z3").Value = 0.76
$ws.Range("AA3").Value = 833.84
$ws.Range("AB3").Value = 162.84
$ws.Range("AC3").Value = 801
$ws.Range("AD3").Value = 8.03
$ws.Range("AE3").Value = 12358
$ws.Range("AF3").Value = 0.52
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 6.22
$ws.Range("AI3").Value = 44.27
$ws.Range("AJ3").Value = 57545890
$ws.Range("U3").ClearContents()

# --- Row 4 ---
$ws.Range("D4").Value = 7616
$ws.Range("E4").Value = 699
$ws.Range("F4").Value = 699
$ws.Range("G4").Value = 678
$ws.Range("H4").Value = 542
$ws.Range("I4").Value = 489
$ws.Range("J4").Value = 53
$ws.Range("K4").Value = 61304
$ws.Range("L4").Value = 53429
$ws.Range("M4").Value = 7875
$ws.Range("N4").Value = 7367
$ws.Range("O4").Value = 508
$ws.Range("P4").Value = 2877
$ws.Range("Q4").Value = 9585
$ws.Range("R4").Value = 243
$ws.Range("S4").Value = -9820
$ws.Range("T4").Value = 33
$ws.Range("V4").Value = 27212
$ws.Range("W4").Value = 9.17
$ws.Range("X4").Value = 7.11
$ws.Range("Y4").Value = 6.77
$ws.Range("Z4").Value = 0.82
$ws.Range("AA4").Value = 678.43
$ws.Range("AB4").Value = 174.39
$ws.Range("AC4").Value = 850
$ws.Range("AD4").Value = 8.2
$ws.Range("AE4").Value = 12866
$ws.Range("AF4").Value = 0.54
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 5.02
$ws.Range("AI4").Value = 30.76
$ws.Range("AJ4").Value = 57545890
$ws.Range("U4").ClearContents()

# --- Row 5 ---
$ws.Range("D5").Value = 6216
$ws.Range("E5").Value = 1126
$ws.Range("F5").Value = 1126
$ws.Range("G5").Value = 739
$ws.Range("H5").Value = 538
$ws.Range("I5").Value = 520
$ws.Range("J5").Value = 18
$ws.Range("K5").Value = 53106
$ws.Range("L5").Value = 46074
$ws.Range("M5").Value = 7032
$ws.Range("N5").Value = 7032
$ws.Range("P5").Value = 2877
$ws.Range("Q5").Value = 4712
$ws.Range("R5").Value = 698
$ws.Range("S5").Value = -5862
$ws.Range("T5").Value = 54
$ws.Range("V5").Value = 27552
$ws.Range("W5").Value = 18.11
$ws.Range("X5").Value = 8.65
$ws.Range("Y5").Value = 7.22
$ws.Range("Z5").Value = 0.94
$ws.Range("AA5").Value = 655.16
$ws.Range("AB5").Value = 145.1
$ws.Range("AC5").Value = 904
$ws.Range("AD5").Value = 7.71
$ws.Range("AE5").Value = 12281
$ws.Range("AF5").Value = 0.57
$ws.Range("AG5").Value = 1653
$ws.Range("AH5").Value = 23.72
$ws.Range("AI5").Value = 177.86
$ws.Range("AJ5").Value = 57545890
$ws.Range("O5").ClearContents()
$ws.Range("U5").ClearContents()

# --- Row 6 ---
$ws.Range("D6").Value = 6152
$ws.Range("E6").Value = 1100
$ws.Range("F6").Value = 1100
$ws.Range("G6").Value = 1099
$ws.Range("H6").Value = 910
$ws.Range("I6").Value = 910
$ws.Range("K6").Value = 62014
$ws.Range("L6").Value = 54491
$ws.Range("M6").Value = 7524
$ws.Range("N6").Value = 7524
$ws.Range("P6").Value = 2877
$ws.Range("Q6").Value = -6908
$ws.Range("R6").Value = -22
$ws.Range("S6").Value = 6433
$ws.Range("T6").Value = 42
$ws.Range("V6").Value = 38299
$ws.Range("W6").Value = 17.88
$ws.Range("X6").Value = 14.79
$ws.Range("Y6").Value = 12.5
$ws.Range("Z6").Value = 1.58
$ws.Range("AA6").Value = 724.26
$ws.Range("AB6").Value = 162.17
$ws.Range("AC6").Value = 1581
$ws.Range("AD6").Value = 5.4
$ws.Range("AE6").Value = 13139
$ws.Range("AF6").Value = 0.65
$ws.Range("AG6").Value = 480
$ws.Range("AH6").Value = 5.63
$ws.Range("AI6").Value = 30.21
$ws.Range("AJ6").Value = 57545890
$ws.Range("U6").ClearContents()

# --- Row 7 ---
$ws.Range("E7").Value = 1270
$ws.Range("G7").Value = 1270
$ws.Range("H7").Value = 970
$ws.Range("I7").Value = 970
$ws.Range("K7").Value = 73400
$ws.Range("L7").Value = 65500
$ws.Range("M7").Value = 7900
$ws.Range("N7").Value = 7900
$ws.Range("P7").Value = 2880
$ws.Range("Y7").Value = 12.58
$ws.Range("Z7").Value = 1.43
$ws.Range("AA7").Value = 829.11
$ws.Range("AC7").Value = 1686
$ws.Range("AD7").Value = 7.59
$ws.Range("AE7").Value = 14023
$ws.Range("AF7").Value = 0.91
$ws.Range("AG7").Value = 580
$ws.Range("AH7").Value = 4.53
$ws.Range("AI7").Value = 34.41
$ws.Range("D7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()

# --- Row 8 ---
$ws.Range("E8").Value = 1360
$ws.Range("G8").Value = 1360
$ws.Range("H8").Value = 1020
$ws.Range("I8").Value = 1020
$ws.Range("K8").Value = 79440
$ws.Range("L8").Value = 70780
$ws.Range("M8").Value = 8660
$ws.Range("N8").Value = 8660
$ws.Range("P8").Value = 2880
$ws.Range("Y8").Value = 12.32
$ws.Range("Z8").Value = 1.34
$ws.Range("AA8").Value = 817.32
$ws.Range("AC8").Value = 1772
$ws.Range("AD8").Value = 7.22
$ws.Range("AE8").Value = 15381
$ws.Range("AF8").Value = 0.83
$ws.Range("AG8").Value = 590
$ws.Range("AH8").Value = 4.61
$ws.Range("AI8").Value = 33.29
$ws.Range("D8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()

# --- Row 9 ---
$ws.Range("E9").Value = 1410
$ws.Range("G9").Value = 1410
$ws.Range("H9").Value = 1060
$ws.Range("I9").Value = 1060
$ws.Range("K9").Value = 87680
$ws.Range("L9").Value = 78230
$ws.Range("M9").Value = 9450
$ws.Range("N9").Value = 9450
$ws.Range("P9").Value = 2880
$ws.Range("Y9").Value = 12.58
$ws.Range("Z9").Value = 1.27
$ws.Range("AA9").Value = 827.83
$ws.Range("AC9").Value = 1842
$ws.Range("AD9").Value = 6.95
$ws.Range("AE9").Value = 16785
$ws.Range("AF9").Value = 0.76
$ws.Range("AG9").Value = 600
$ws.Range("AH9").Value = 4.69
$ws.Range("AI9").Value = 32.57
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
